$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 456350
$ws.Range("D2").Value = 677856754
$ws.Range("C9").Value = 1508
$ws.Range("D9").Value = 4141244
$ws.Range("C11").Value = 177423
$ws.Range("D11").Value = 436090543
$ws.Range("C14").Value = 92265
$ws.Range("D14").Value = 216652328
$ws.Range("C18").Value = 4957
$ws.Range("D18").Value = 8446401
$ws.Range("C21").Value = 97
$ws.Range("D21").Value = 223464
$ws.Range("C23").Value = 11846
$ws.Range("D23").Value = 24379328
$ws.Range("C25").Value = 108660
$ws.Range("D25").Value = 158592060
$ws.Range("C31").Value = 47972
$ws.Range("D31").Value = 114528453
$ws.Range("C34").Value = 17421
$ws.Range("D34").Value = 40227622
$ws.Range("C37").Value = 2001
$ws.Range("D37").Value = 3540829
$ws.Range("C39").Value = 3038
$ws.Range("D39").Value = 5940586
$ws.Range("C40").Value = 133837
$ws.Range("D40").Value = 198856387
$ws.Range("C48").Value = 67221
$ws.Range("D48").Value = 165768659
$ws.Range("C51").Value = 14277
$ws.Range("D51").Value = 34107664
$ws.Range("C53").Value = 1730
$ws.Range("D53").Value = 2844397
$ws.Range("C56").Value = 3822
$ws.Range("D56").Value = 7588638
$ws.Range("C57").Value = 96247
$ws.Range("D57").Value = 142767023
$ws.Range("C65").Value = 41515
$ws.Range("D65").Value = 96999551
$ws.Range("C68").Value = 16633
$ws.Range("D68").Value = 37401413
$ws.Range("C70").Value = 1768
$ws.Range("D70").Value = 3167178
$ws.Range("C74").Value = 2502
$ws.Range("D74").Value = 4972370
$ws.Range("C76").Value = 28646
$ws.Range("D76").Value = 44149632
$ws.Range("C80").Value = 11862
$ws.Range("D80").Value = 31414048
$ws.Range("C82").Value = 8023
$ws.Range("D82").Value = 19503578
$ws.Range("C85").Value = 499
$ws.Range("D85").Value = 983098
$ws.Range("C86").Value = 200153
$ws.Range("D86").Value = 301579474
$ws.Range("C87").Value = 141
$ws.Range("D87").Value = 418519
$ws.Range("C92").Value = 94142
$ws.Range("D92").Value = 221516178
$ws.Range("C95").Value = 44534
$ws.Range("D95").Value = 101452882
$ws.Range("C98").Value = 6945
$ws.Range("D98").Value = 26295361
$ws.Range("C100").Value = 4682
$ws.Range("D100").Value = 9384119
$ws.Range("C102").Value = 44478
$ws.Range("D102").Value = 65968613
$ws.Range("C106").Value = 10830
$ws.Range("D106").Value = 18768471
$ws.Range("C108").Value = 10136
$ws.Range("D108").Value = 16858490
$ws.Range("C111").Value = 663
$ws.Range("D111").Value = 1074746
$ws.Range("C112").Value = 18570
$ws.Range("D112").Value = 37898062
$ws.Range("C115").Value = 4176
$ws.Range("D115").Value = 9104278
$ws.Range("C117").Value = 6057
$ws.Range("D117").Value = 13643156
$ws.Range("C122").Value = 200296
$ws.Range("D122").Value = 291366130
$ws.Range("C127").Value = 1610
$ws.Range("D127").Value = 4457750
$ws.Range("C129").Value = 78315
$ws.Range("D129").Value = 184915781
$ws.Range("C132").Value = 41993
$ws.Range("D132").Value = 95563597
$ws.Range("C135").Value = 1802
$ws.Range("D135").Value = 3674301
$ws.Range("C139").Value = 3894
$ws.Range("D139").Value = 7920845
$ws.Range("C141").Value = 801177
$ws.Range("D141").Value = 1276930730
$ws.Range("C142").Value = 140
$ws.Range("D142").Value = 250663
$ws.Range("C146").Value = 2782
$ws.Range("D146").Value = 9619112
$ws.Range("C148").Value = 323956
$ws.Range("D148").Value = 776604998
$ws.Range("C149").Value = 1072
$ws.Range("D149").Value = 4647241
$ws.Range("C151").Value = 300789
$ws.Range("D151").Value = 677141237
$ws.Range("C154").Value = 3807
$ws.Range("D154").Value = 6409069
$ws.Range("C157").Value = 11864
$ws.Range("D157").Value = 24680597
$ws.Range("C160").Value = 54776
$ws.Range("D160").Value = 76109114
$ws.Range("C166").Value = 16487
$ws.Range("D166").Value = 26617558
$ws.Range("C168").Value = 4563
$ws.Range("D168").Value = 7233389
$ws.Range("C174").Value = 24626
$ws.Range("D174").Value = 36710816
$ws.Range("C178").Value = 10359
$ws.Range("D178").Value = 21141376
$ws.Range("C180").Value = 7295
$ws.Range("D180").Value = 13785766
$ws.Range("C183").Value = 435
$ws.Range("D183").Value = 795867
$ws.Range("C185").Value = 33336
$ws.Range("D185").Value = 80996378
$ws.Range("C187").Value = 3259
$ws.Range("D187").Value = 7642574
$ws.Range("C192").Value = 123484
$ws.Range("D192").Value = 184886253
$ws.Range("C198").Value = 1163
$ws.Range("D198").Value = 3488298
$ws.Range("C200").Value = 51230
$ws.Range("D200").Value = 125769618
$ws.Range("C201").Value = 57
$ws.Range("D201").Value = 196545
$ws.Range("C202").Value = 20226
$ws.Range("D202").Value = 47575830
$ws.Range("C204").Value = 1590
$ws.Range("D204").Value = 2736019
$ws.Range("C207").Value = 2820
$ws.Range("D207").Value = 5750340
$ws.Range("C209").Value = 326777
$ws.Range("D209").Value = 470052629
$ws.Range("C218").Value = 126742
$ws.Range("D218").Value = 296742418
$ws.Range("C221").Value = 49874
$ws.Range("D221").Value = 113021934
$ws.Range("C224").Value = 6551
$ws.Range("D224").Value = 11917591
$ws.Range("C227").Value = 8152
$ws.Range("D227").Value = 15562073
$ws.Range("C230").Value = 371016
$ws.Range("D230").Value = 524585758
$ws.Range("C231").Value = 205
$ws.Range("D231").Value = 312823
$ws.Range("C239").Value = 141065
$ws.Range("D239").Value = 329543755
$ws.Range("C242").Value = 78010
$ws.Range("D242").Value = 175722191
$ws.Range("C245").Value = 5992
$ws.Range("D245").Value = 9920235
$ws.Range("C248").Value = 45
$ws.Range("D248").Value = 179640
$ws.Range("C249").Value = 10288
$ws.Range("D249").Value = 19854337
$ws.Range("C252").Value = 147757
$ws.Range("D252").Value = 217919617
$ws.Range("C257").Value = 1035
$ws.Range("D257").Value = 2986745
$ws.Range("C259").Value = 73366
$ws.Range("D259").Value = 175252784
$ws.Range("C262").Value = 19317
$ws.Range("D262").Value = 44189892
$ws.Range("C264").Value = 2372
$ws.Range("D264").Value = 4271392
$ws.Range("C266").Value = 4057
$ws.Range("D266").Value = 8098295
$ws.Range("C267").Value = 371616
$ws.Range("D267").Value = 542505211
$ws.Range("C268").Value = 242
$ws.Range("D268").Value = 389001
$ws.Range("C274").Value = 1585
$ws.Range("D274").Value = 4942000
$ws.Range("C276").Value = 147111
$ws.Range("D276").Value = 357189169
$ws.Range("C279").Value = 104113
$ws.Range("D279").Value = 241943196
$ws.Range("C282").Value = 3265
$ws.Range("D282").Value = 5629994
$ws.Range("C285").Value = 8396
$ws.Range("D285").Value = 16973801
